$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 25
$n = 24

$colB = New-Object 'object[,]' $n,1
$colB[0,0] = 0.3483199873666081
$colB[1,0] = 0.3041225684244182
$colB[2,0] = 0.2768875659509717
$colB[3,0] = 0.2657652826387391
$colB[4,0] = 0.2639170206162476
$colB[5,0] = 0.276737662317089
$colB[6,0] = 0.333101408955315
$colB[7,0] = 0.4428285929688229
$colB[8,0] = 0.5229271360309156
$colB[9,0] = 0.5592479070427316
$colB[10,0] = 0.5729842494943966
$colB[11,0] = 0.5700266753715084
$colB[12,0] = 0.5603783603944521
$colB[13,0] = 0.5544661803445194
$colB[14,0] = 0.520551076407969
$colB[15,0] = 0.4997148611580542
$colB[16,0] = 0.4877195071149458
$colB[17,0] = 0.4836562362477821
$colB[18,0] = 0.501934047213382
$colB[19,0] = 0.5632127868976795
$colB[20,0] = 0.6031593510644768
$colB[21,0] = 0.5818487717590983
$colB[22,0] = 0.5009308040235112
$colB[23,0] = 0.4132333661295036
$ws.Range("B" + $startRow + ":B" + $endRow).Value = $colB

$colC = New-Object 'object[,]' $n,1
$colC[0,0] = 0.07766027381389051
$colC[1,0] = 0.07123906844635997
$colC[2,0] = 0.0672741386793092
$colC[3,0] = 0.06565288561210991
$colC[4,0] = 0.06538334757070174
$colC[5,0] = 0.06725229609048711
$colC[6,0] = 0.07545091995275754
$colC[7,0] = 0.09134832867559339
$colC[8,0] = 0.1029151245213171
$colC[9,0] = 0.1081520041497157
$colC[10,0] = 0.1101314162769
$colC[11,0] = 0.1097052797036895
$colC[12,0] = 0.1083149258977727
$colC[13,0] = 0.1074628118397243
$colC[14,0] = 0.1025723724002319
$colC[15,0] = 0.09956579522084041
$colC[16,0] = 0.0978341523478008
$colC[17,0] = 0.09724744906297644
$colC[18,0] = 0.09988609330538623
$colC[19,0] = 0.1087234071785019
$colC[20,0] = 0.1144776048203084
$colC[21,0] = 0.1114084804535622
$colC[22,0] = 0.09974129625958028
$colC[23,0] = 0.08706727558663374
$ws.Range("C" + $startRow + ":C" + $endRow).Value = $colC

$colE = New-Object 'object[,]' $n,1
$colE[0,0] = 0.4214918725342187
$colE[1,0] = 0.3678202120443927
$colE[2,0] = 0.3349375038587254
$colE[3,0] = 0.3215542924570798
$colE[4,0] = 0.3193329940751966
$colE[5,0] = 0.3347569470696925
$colE[6,0] = 0.4029701469637672
$colE[7,0] = 0.5373724907827011
$colE[8,0] = 0.6366061295298806
$colE[9,0] = 0.6818782999108919
$colE[10,0] = 0.6990420370345873
$colE[11,0] = 0.695344600418224
$colE[12,0] = 0.6832899589519599
$colE[13,0] = 0.6759088090582708
$colE[14,0] = 0.6336502290723445
$colE[15,0] = 0.6077603318098141
$colE[16,0] = 0.5928813261688362
$colE[17,0] = 0.5878455992767186
$colE[18,0] = 0.6105150849500518
$colE[19,0] = 0.6868301412168307
$colE[20,0] = 0.7368245474456927
$colE[21,0] = 0.7101303157136414
$colE[22,0] = 0.6092696442384948
$colE[23,0] = 0.5009346304369728
$ws.Range("E" + $startRow + ":E" + $endRow).Value = $colE

$colF = New-Object 'object[,]' $n,1
$colF[0,0] = 0.4443680307746121
$colF[1,0] = 0.3878228170618172
$colF[2,0] = 0.3531389305169483
$colF[3,0] = 0.3390132514313251
$colF[4,0] = 0.336668177824194
$colF[5,0] = 0.3529483938344953
$colF[6,0] = 0.4248636149813478
$colF[7,0] = 0.5661985755041457
$colF[8,0] = 0.6702781546542269
$colF[9,0] = 0.7176906081379002
$colF[10,0] = 0.7356546913071611
$colF[11,0] = 0.7317853510981394
$colF[12,0] = 0.7191683204515869
$colF[13,0] = 0.7114413442032514
$colF[14,0] = 0.6671810134426437
$colF[15,0] = 0.6400460337125793
$colF[16,0] = 0.6244449056556647
$colF[17,0] = 0.619163680173358
$colF[18,0] = 0.642933953830422
$colF[19,0] = 0.7228739723491628
$colF[20,0] = 0.7751780083420101
$colF[21,0] = 0.7472568307830727
$colF[22,0] = 0.6416283278902171
$colF[23,0] = 0.5279251897347166
$ws.Range("F" + $startRow + ":F" + $endRow).Value = $colF

$colG = New-Object 'object[,]' $n,1
$colG[0,0] = 0.2517007221420897
$colG[1,0] = 0.2562701936192191
$colG[2,0] = 0.2593547557453135
$colG[3,0] = 0.260681616736008
$colG[4,0] = 0.260906154002523
$colG[5,0] = 0.2593723676623974
$colG[6,0] = 0.2532182415693995
$colG[7,0] = 0.2433741370383586
$colG[8,0] = 0.2375130128124638
$colG[9,0] = 0.2351477268900055
$colG[10,0] = 0.2342956078055352
$colG[11,0] = 0.2344771848777185
$colG[12,0] = 0.2350767478375104
$colG[13,0] = 0.2354496791349163
$colG[14,0] = 0.2376736703893272
$colG[15,0] = 0.2391153047608228
$colG[16,0] = 0.2399728063451363
$colG[17,0] = 0.2402679955290736
$colG[18,0] = 0.2389589081123518
$colG[19,0] = 0.2348994570996226
$colG[20,0] = 0.2325004236523682
$colG[21,0] = 0.2337574932162525
$colG[22,0] = 0.2390295257001966
$colG[23,0] = 0.2457975228181866
$ws.Range("G" + $startRow + ":G" + $endRow).Value = $colG

$colH = New-Object 'object[,]' $n,1
$colH[0,0] = 0.4426920784174513
$colH[1,0] = 0.4485126668844259
$colH[2,0] = 0.4523321796906146
$colH[3,0] = 0.4539503894675541
$colH[4,0] = 0.4542228188944968
$colH[5,0] = 0.4523537535601143
$colH[6,0] = 0.4446480077171771
$colH[7,0] = 0.431487970616601
$colH[8,0] = 0.4230109194674512
$colH[9,0] = 0.4194137397873661
$colH[10,0] = 0.4180888861583938
$colH[11,0] = 0.4183725562111391
$colH[12,0] = 0.419303994950539
$colH[13,0] = 0.4198793902307614
$colH[14,0] = 0.4232512224542546
$colH[15,0] = 0.4253861339142588
$colH[16,0] = 0.4266384597699258
$colH[17,0] = 0.4270666612937504
$colH[18,0] = 0.4251563450999498
$colH[19,0] = 0.4190293954269819
$colH[20,0] = 0.4152426414250456
$colH[21,0] = 0.4172437745329773
$colH[22,0] = 0.4252601549488588
$colH[23,0] = 0.434839015018369
$ws.Range("H" + $startRow + ":H" + $endRow).Value = $colH

$colI = New-Object 'object[,]' $n,1
$colI[0,0] = 0.3261503604462632
$colI[1,0] = 0.3326743421453919
$colI[2,0] = 0.3369324167165537
$colI[3,0] = 0.3387309529775084
$colI[4,0] = 0.339033420615225
$colI[5,0] = 0.3369564161091212
$colI[6,0] = 0.3283473927233533
$colI[7,0] = 0.3134719984428695
$colI[8,0] = 0.3037727442139211
$colI[9,0] = 0.2996286515988071
$colI[10,0] = 0.298098071865498
$colI[11,0] = 0.2984259862681498
$colI[12,0] = 0.2995019535165468
$colI[13,0] = 0.3001660583009897
$colI[14,0] = 0.3040489780910764
$colI[15,0] = 0.3064998013017934
$colI[16,0] = 0.3079346809102166
$colI[17,0] = 0.3084248356940815
$colI[18,0] = 0.3062362946741306
$colI[19,0] = 0.2991848644869766
$colI[20,0] = 0.2948019746010981
$colI[21,0] = 0.2971205145226961
$colI[22,0] = 0.306355345450255
$colI[23,0] = 0.3172806915296249
$ws.Range("I" + $startRow + ":I" + $endRow).Value = $colI

$colK = New-Object 'object[,]' $n,1
$colK[0,0] = 0.3790457237252554
$colK[1,0] = 0.3313300743348293
$colK[2,0] = 0.3019030122952699
$colK[3,0] = 0.2898794577183423
$colK[4,0] = 0.287881058365997
$colK[5,0] = 0.3017409861228941
$colK[6,0] = 0.3626206814189459
$colK[7,0] = 0.4809503835636235
$colK[8,0] = 0.5672146765993773
$colK[9,0] = 0.6063068645489977
$colK[10,0] = 0.6210878540784392
$colK[11,0] = 0.6179055116291465
$colK[12,0] = 0.6075233585464446
$colK[13,0] = 0.6011610519791191
$colK[14,0] = 0.5646568239732517
$colK[15,0] = 0.5422236900637927
$colK[16,0] = 0.5293066869759855
$colK[17,0] = 0.5249308223734772
$colK[18,0] = 0.5446131963666403
$colK[19,0] = 0.610573462316097
$colK[20,0] = 0.6535515204639353
$colK[21,0] = 0.6306255595343657
$colK[22,0] = 0.5435329624786789
$colK[23,0] = 0.4490550746224642
$ws.Range("K" + $startRow + ":K" + $endRow).Value = $colK

$colO = New-Object 'object[,]' $n,1
$colO[0,0] = 1.307428150239318
$colO[1,0] = 1.329078728740186
$colO[2,0] = 1.343470494376788
$colO[3,0] = 1.349610784280223
$colO[4,0] = 1.350646997576206
$colO[5,0] = 1.3435521897065
$colO[6,0] = 1.314664957135392
$colO[7,0] = 1.266758966179779
$colO[8,0] = 1.236929569081298
$colO[9,0] = 1.224532973981908
$colO[10,0] = 1.220008065970717
$colO[11,0] = 1.220975039289456
$colO[12,0] = 1.224157307507951
$colO[13,0] = 1.226128623439607
$colO[14,0] = 1.237763382107687
$colO[15,0] = 1.245201879671299
$colO[16,0] = 1.249590665488839
$colO[17,0] = 1.251095564490768
$colO[18,0] = 1.244398611788327
$colO[19,0] = 1.223217994418036
$colO[20,0] = 1.210363054279526
$colO[21,0] = 1.217133346063179
$colO[22,0] = 1.244761419457191
$colO[23,0] = 1.278779022084677
$ws.Range("O" + $startRow + ":O" + $endRow).Value = $colO
